# Updates crypto price (D) and 1h-volume-change (E) columns to the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value.
$updates = @{
    'D2' = '43.690.92'
    'E2' = '  -0.37%  '
    'D3' = '2.282.69'
    'E3' = '  -0.46%  '
    'D4' = '1.00'
    'E4' = '  +0.02%  '
    'D5' = '113.75'
    'E5' = '  +10.62%  '
    'D6' = '267.78'
    'E6' = '  -1.06%  '
    'E7' = '  +0.21%  '
    'E8' = '  +0.11%  '
    'D9' = '0.613'
    'E9' = '  +0.95%  '
    'D10' = '48.60'
    'E10' = '  +4.86%  '
    'E11' = '  +0.37%  '
    'E12' = '  +9.66%  '
    'E13' = '  +0.68%  '
    'D14' = '15.78'
    'E14' = '  +1.50%  '
    'D15' = '2.624.13'
    'E15' = '  -0.29%  '
    'D16' = '0.877'
    'E16' = '  +2.31%  '
    'D17' = '2.284.18'
    'E17' = '  -0.11%  '
    'D18' = '43.527.00'
    'E18' = '  -0.49%  '
    'E19' = '  -1.18%  '
    'D20' = '7.06'
    'E20' = '  +12.67%  '
    'D21' = '72.08'
    'E21' = '  -0.32%  '
    'E22' = '  -3.75%  '
    'D23' = '9.95'
    'E23' = '  +7.27%  '
    'D24' = '232.84'
    'E24' = '  -0.26%  '
    'D25' = '2.88'
    'E25' = '  +0.18%  '
    'E26' = '  -0.01%  '
    'D27' = '11.58'
    'E27' = '  +3.01%  '
    'D28' = '41.53'
    'E28' = '  +0.15%  '
    'E29' = '  -1.58%  '
    'E30' = '  -1.49%  '
    'D31' = '173.28'
    'E31' = '  -2.43%  '
    'E32' = '  -1.31%  '
    'D33' = '0.0908'
    'E33' = '  +0.61%  '
    'D34' = '5.68'
    'E34' = '  +3.18%  '
    'E35' = '  +0.24%  '
    'D36' = '4.64'
    'E36' = '  -4.52%  '
    'E37' = '  -1.81%  '
    'E38' = '  -3.68%  '
    'E39' = '  +4.97%  '
    'D40' = '14.69'
    'E40' = '  +19.98%  '
    'D41' = '75.24'
    'E41' = '  +14.11%  '
    'D42' = '2.43'
    'E42' = '  +4.31%  '
    'E43' = '  +0.39%  '
    'D44' = '6.26'
    'E44' = '  +18.81%  '
    'E45' = '  +0.05%  '
    'E46' = '  +0.19%  '
    'D47' = '8.69'
    'E47' = '  -1.32%  '
    'E48' = '  +1.86%  '
    'D49' = '102.04'
    'E49' = '  +2.92%  '
    'D50' = '0.0998'
    'E50' = '  -2.04%  '
    'D51' = '0.455'
    'E51' = '  +2.53%  '
}

foreach ($ref in $updates.Keys) {
    $newValue = $updates[$ref]
    $cell = $ws.Range($ref)

    # Columns D values are plain-looking numbers stored as text in the sheet (e.g. "1.00",
    # "48.60", "43.690.92"). Assigning such a string directly would make Excel silently
    # reinterpret/renormalize it as a real number (losing the exact text, e.g. "1.00" -> 1).
    # Force the cell to Text format first so the exact string is preserved, then restore
    # the cells original style so no visible formatting changes are introduced.
    if ($ref.StartsWith("D")) {
        $origStyle = $cell.Style
        $cell.NumberFormat = "@"
        $cell.Value = $newValue
        $cell.Style = $origStyle
    } else {
        $cell.Value = $newValue
    }
}
